# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Cells.Item(2, 4).Value = "'43.615.67"
$ws.Cells.Item(2, 5).Value = "'  -0.85%  "

# Row 3: 'Ethereum'
$ws.Cells.Item(3, 4).Value = "'2.224.61"
$ws.Cells.Item(3, 5).Value = "'  -1.81%  "

# Row 4: 'TetherUSD'
$ws.Cells.Item(4, 5).Value = "'  +0.25%  "

# Row 5: 'BNB'
$ws.Cells.Item(5, 4).Value = "'313.17"
$ws.Cells.Item(5, 5).Value = "'  -1.87%  "

# Row 6: 'Solana'
$ws.Cells.Item(6, 4).Value = "'98.04"
$ws.Cells.Item(6, 5).Value = "'  -4.45%  "

# Row 7: 'XRP'
$ws.Cells.Item(7, 5).Value = "'  -3.05%  "

# Row 8: 'USDC'
$ws.Cells.Item(8, 5).Value = "'  +0.14%  "

# Row 9: 'Cardano'
$ws.Cells.Item(9, 4).Value = "'0.534"
$ws.Cells.Item(9, 5).Value = "'  -6.46%  "

# Row 10: 'Avalanche'
$ws.Cells.Item(10, 4).Value = "'36.03"
$ws.Cells.Item(10, 5).Value = "'  -6.02%  "

# Row 11: 'Dogecoin'
$ws.Cells.Item(11, 4).Value = "'0.0821"
$ws.Cells.Item(11, 5).Value = "'  -2.02%  "

# Row 12: 'Polkadot'
$ws.Cells.Item(12, 4).Value = "'7.38"
$ws.Cells.Item(12, 5).Value = "'  -6.07%  "

# Row 13: 'TRON'
$ws.Cells.Item(13, 5).Value = "'  -3.20%  "

# Row 14: 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 4).Value = "'2.563.45"
$ws.Cells.Item(14, 5).Value = "'  -1.78%  "

# Row 15: 'WrappedEther' -> 'Polygon'
$ws.Cells.Item(15, 2).Value = "'Polygon"
$ws.Cells.Item(15, 3).Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(15, 4).Value = "'0.838"
$ws.Cells.Item(15, 5).Value = "'  -3.95%  "

# Row 16: 'Polygon' -> 'Chainlink'
$ws.Cells.Item(16, 2).Value = "'Chainlink"
$ws.Cells.Item(16, 3).Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(16, 4).Value = "'14.08"
$ws.Cells.Item(16, 5).Value = "'  -3.24%  "

# Row 17: 'Chainlink' -> 'WrappedEther'
$ws.Cells.Item(17, 2).Value = "'WrappedEther"
$ws.Cells.Item(17, 3).Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(17, 4).Value = "'2.223.05"
$ws.Cells.Item(17, 5).Value = "'  -2.24%  "

# Row 18: 'WrappedBTC'
$ws.Cells.Item(18, 4).Value = "'43.478.87"
$ws.Cells.Item(18, 5).Value = "'  -0.91%  "

# Row 19: 'InternetComputer(DFINITY)'
$ws.Cells.Item(19, 4).Value = "'13.03"
$ws.Cells.Item(19, 5).Value = "'  -9.95%  "

# Row 20: 'ShibaInu'
$ws.Cells.Item(20, 4).Value = "'0.0₃0963"
$ws.Cells.Item(20, 5).Value = "'  -2.96%  "

# Row 21: 'Uniswap'
$ws.Cells.Item(21, 4).Value = "'6.31"
$ws.Cells.Item(21, 5).Value = "'  -5.29%  "

# Row 22: 'Litecoin'
$ws.Cells.Item(22, 4).Value = "'65.24"
$ws.Cells.Item(22, 5).Value = "'  -1.21%  "

# Row 23: 'BitcoinCash'
$ws.Cells.Item(23, 4).Value = "'234.80"
$ws.Cells.Item(23, 5).Value = "'  -1.78%  "

# Row 24: 'PancakeSwap'
$ws.Cells.Item(24, 5).Value = "'  -7.07%  "

# Row 25: 'ImmutableX'
$ws.Cells.Item(25, 5).Value = "'  -7.24%  "

# Row 26: 'Dai'
$ws.Cells.Item(26, 5).Value = "'  +0.15%  "

# Row 27: 'Cosmos'
$ws.Cells.Item(27, 4).Value = "'10.00"
$ws.Cells.Item(27, 5).Value = "'  -2.39%  "

# Row 28: 'Toncoin'
$ws.Cells.Item(28, 4).Value = "'2.23"
$ws.Cells.Item(28, 5).Value = "'  +1.31%  "

# Row 29: 'InjectiveProtocol'
$ws.Cells.Item(29, 4).Value = "'36.19"
$ws.Cells.Item(29, 5).Value = "'  -6.12%  "

# Row 30: 'Monero' -> 'Filecoin'
$ws.Cells.Item(30, 2).Value = "'Filecoin"
$ws.Cells.Item(30, 3).Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30, 4).Value = "'5.95"
$ws.Cells.Item(30, 5).Value = "'  -8.99%  "

# Row 31: 'Filecoin' -> 'Monero'
$ws.Cells.Item(31, 2).Value = "'Monero"
$ws.Cells.Item(31, 3).Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(31, 4).Value = "'159.30"
$ws.Cells.Item(31, 5).Value = "'  -2.50%  "

# Row 32: 'EthereumClassic'
$ws.Cells.Item(32, 4).Value = "'19.89"
$ws.Cells.Item(32, 5).Value = "'  -2.88%  "

# Row 33: 'Hedera'
$ws.Cells.Item(33, 4).Value = "'0.0829"
$ws.Cells.Item(33, 5).Value = "'  -6.19%  "

# Row 34: 'WEMIXToken'
$ws.Cells.Item(34, 4).Value = "'2.68"
$ws.Cells.Item(34, 5).Value = "'  -1.20%  "

# Row 35: 'LidoDAOToken'
$ws.Cells.Item(35, 4).Value = "'3.17"
$ws.Cells.Item(35, 5).Value = "'  -2.83%  "

# Row 36: 'ARBITRUM' -> 'Kaspa'
$ws.Cells.Item(36, 2).Value = "'Kaspa"
$ws.Cells.Item(36, 3).Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(36, 4).Value = "'0.109"
$ws.Cells.Item(36, 5).Value = "'  +1.24%  "

# Row 37: 'Kaspa' -> 'ARBITRUM'
$ws.Cells.Item(37, 2).Value = "'ARBITRUM"
$ws.Cells.Item(37, 3).Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(37, 4).Value = "'1.89"
$ws.Cells.Item(37, 5).Value = "'  -6.47%  "

# Row 38: 'Stellar'
$ws.Cells.Item(38, 5).Value = "'  -3.48%  "

# Row 39: 'Celestia'
$ws.Cells.Item(39, 4).Value = "'15.57"
$ws.Cells.Item(39, 5).Value = "'  -1.99%  "

# Row 40: 'NEARProtocol'
$ws.Cells.Item(40, 4).Value = "'3.57"
$ws.Cells.Item(40, 5).Value = "'  -7.80%  "

# Row 41: 'RenderToken'
$ws.Cells.Item(41, 4).Value = "'4.00"
$ws.Cells.Item(41, 5).Value = "'  -12.15%  "

# Row 42: 'VeChain'
$ws.Cells.Item(42, 5).Value = "'  -5.76%  "

# Row 43: 'FirstDigitalUSD'
$ws.Cells.Item(43, 5).Value = "'  +0.23%  "

# Row 44: 'Maker'
$ws.Cells.Item(44, 4).Value = "'1.704.29"
$ws.Cells.Item(44, 5).Value = "'  -4.18%  "

# Row 45: 'BitcoinSV'
$ws.Cells.Item(45, 4).Value = "'82.06"
$ws.Cells.Item(45, 5).Value = "'  -3.50%  "

# Row 46: 'Algorand'
$ws.Cells.Item(46, 4).Value = "'0.195"
$ws.Cells.Item(46, 5).Value = "'  -5.92%  "

# Row 47: 'THORChain'
$ws.Cells.Item(47, 5).Value = "'  -5.46%  "

# Row 48: 'ordi'
$ws.Cells.Item(48, 4).Value = "'72.65"
$ws.Cells.Item(48, 5).Value = "'  -2.47%  "

# Row 49: 'Aave'
$ws.Cells.Item(49, 4).Value = "'101.63"
$ws.Cells.Item(49, 5).Value = "'  -2.87%  "

# Row 50: 'Stacks'
$ws.Cells.Item(50, 5).Value = "'  +1.28%  "

# Row 51: 'MultiversX'
$ws.Cells.Item(51, 4).Value = "'56.45"
$ws.Cells.Item(51, 5).Value = "'  -5.59%  "
